$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item Stack row: base cost / grow rate pairs for the new 101 entry (row 5)
$ws.Range("A5").Value = 101
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1.5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1.3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.2
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1.2
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1.2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.2

# The newly filled-in cells (E5:M5) need the same formatting as the rest of
# the row (A5:D5 already carry it) - vertical-centered, wrapped text.
$ws.Range("E5:M5").WrapText = $true
$ws.Range("E5:M5").VerticalAlignment = -4108

# Resize the (now wider) data columns to fit their contents again
# (column A and H are left at the default width, same as the source sheet)
$ws.Columns.Item(2).ColumnWidth = 19.65
$ws.Columns.Item(3).ColumnWidth = 20.5
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 15.36
$ws.Columns.Item(6).ColumnWidth = 12.5
$ws.Columns.Item(7).ColumnWidth = 13.22
$ws.Columns.Item(9).ColumnWidth = 24.79
$ws.Columns.Item(10).ColumnWidth = 20.22
$ws.Columns.Item(11).ColumnWidth = 20.93
$ws.Columns.Item(12).ColumnWidth = 23.79
$ws.Columns.Item(13).ColumnWidth = 24.65

# Leave the selection where the user ended up editing
$ws.Range("J13").Select()
